# Apply the "fine tuned EE, ME, and added Transportation Engineering" edit
# to the ME_Programs workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program_choosing")

# --- Update existing "Choose" decisions ---
# RWTH Aachen AUTO: No -> Yes
$ws.Range("B2").Value = "Yes"

# KIT_ME: Yes -> No
$ws.Range("B12").Value = "No"

# --- Add a new "German" (德語) note column for programs requiring German ---
$ws.Range("C3").Value = "德語"   # TUM_MW
$ws.Range("C9").Value = "德語"   # TUBerlin_ME
$ws.Range("C10").Value = "德語"  # RWTH Aachen_ME
$ws.Range("C11").Value = "德語"  # TUBraunschweig_ME

# --- Update selection to match the saved view ---
$ws.Range("A9").Select()
